# Apply the commit's edits to the workbook:
# 1. Rename "Requested quantity" header -> "Weekly_PO_Qty" on "Weekly Quantity" sheet (B1)
# 2. Rename "Requested quantity" header -> "Monthly_PO_Qty" on "Monthly Trend" sheet (B1)
# 3. Add a new "PO Forecast" worksheet (after "Monthly Trend") with ds/PO_Forecast/yhat_lower/yhat_upper data

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- 1 & 2: header renames ---
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3: new "PO Forecast" sheet, placed after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match sheet-level properties (outline + page margins) used on the other sheets
$wsForecast.Outline.SummaryRow = 1
$wsForecast.Outline.SummaryColumn = 1
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match header formatting (bold + border + centered) used on the other sheets
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$rows = @(
    @(44934.99999999999, 125, -105.4592128923455, 351.9174648135174),
    @(44941.99999999999, 126, -95.12906151079086, 342.2976826802428),
    @(44955.99999999999, 126, -105.0403290320058, 345.6976175804758),
    @(44962.99999999999, 127, -84.37971336992382, 360.6152622669971),
    @(44969.99999999999, 127, -113.6191637163692, 360.9335697239458),
    @(44976.99999999999, 127, -92.1067422080724, 339.2011432125562),
    @(44983.99999999999, 128, -103.8464182295613, 361.0518096717302),
    @(44990.99999999999, 128, -117.3892671513755, 352.6937040672702),
    @(44997.99999999999, 129, -106.1967550606527, 373.3831596645383),
    @(45004.99999999999, 129, -101.7807863272088, 362.9878039422624),
    @(45165.99999999999, 137, -89.67679751694551, 345.6600499721617),
    @(45179.99999999999, 138, -91.46237736880194, 384.5022313337099),
    @(45186.99999999999, 139, -84.32247613419759, 358.7846510906555),
    @(45193.99999999999, 139, -90.84480578692907, 372.9522420154984),
    @(45200.99999999999, 139, -89.00993142845812, 359.920997618648),
    @(45207.99999999999, 140, -90.53316187355243, 355.7818418743908),
    @(45214.99999999999, 140, -90.56108328337446, 364.6980385732058),
    @(45221.99999999999, 140, -95.847968469782, 367.0199683439031),
    @(45228.99999999999, 141, -71.22384515749381, 363.4327917680287),
    @(45235.99999999999, 141, -80.59488135768714, 373.2711234344823),
    @(45242.99999999999, 142, -81.30632858681238, 371.8123606388629),
    @(45249.99999999999, 142, -97.34268713635853, 362.6248534650116),
    @(45256.99999999999, 142, -100.2611676291078, 374.5670241529363),
    @(45263.99999999999, 143, -78.80451091241494, 396.0502178884693),
    @(45270.99999999999, 143, -102.7126208149472, 378.7802665961119),
    @(45277.99999999999, 143, -81.80518887199358, 386.8430175886637),
    @(45298.99999999999, 145, -68.54447686283127, 374.5859911548584),
    @(45305.99999999999, 145, -95.12244695301725, 362.9953125566692),
    @(45326.99999999999, 146, -65.6123097415045, 378.9992121181263),
    @(45333.99999999999, 146, -76.59936498666315, 373.7792213909237),
    @(45340.99999999999, 147, -77.62536215378427, 373.2532577366288),
    @(45347.99999999999, 147, -88.65958323878623, 377.7431011310586),
    @(45354.99999999999, 147, -86.51921874671862, 373.926709957817),
    @(45361.99999999999, 148, -92.52482069519348, 371.623527627217),
    @(45368.99999999999, 148, -70.32790516322851, 378.4404716060225),
    @(45375.99999999999, 149, -74.77958571830109, 387.4821473047028),
    @(45382.99999999999, 149, -81.86361931572439, 364.8490154722206),
    @(45389.99999999999, 149, -79.80381091439799, 377.9798084035903),
    @(45396.99999999999, 150, -83.24814662731232, 384.3158432274022),
    @(45403.99999999999, 150, -90.59629394555424, 359.1923409674858),
    @(45410.99999999999, 150, -79.73510147247153, 381.2464778449815),
    @(45417.99999999999, 151, -85.58206228709399, 388.7532749676544),
    @(45424.99999999999, 151, -59.95509580648088, 383.3195811331963),
    @(45431.99999999999, 152, -74.23499232976349, 383.5103324931361),
    @(45438.99999999999, 152, -67.4542698700667, 392.9819555604608),
    @(45452.99999999999, 153, -87.61796905807117, 366.2803850551072),
    @(45459.99999999999, 153, -76.45077894328875, 364.130584124212),
    @(45473.99999999999, 154, -55.16146353430112, 380.1107633603912),
    @(45480.99999999999, 154, -70.47957703015868, 391.3581112495831),
    @(45487.99999999999, 155, -77.43573425549353, 392.590726716648),
    @(45494.99999999999, 155, -79.27287783552715, 390.0958235123881),
    @(45501.99999999999, 155, -77.78233669060548, 367.6901470321844),
    @(45508.99999999999, 156, -75.04382596232981, 375.5435884587087),
    @(45515.99999999999, 156, -60.92234285668302, 381.6946714610437),
    @(45522.99999999999, 156, -76.7494381547759, 373.9238041459069),
    @(45529.99999999999, 157, -89.20500895884156, 384.4940158616038),
    @(45536.99999999999, 157, -74.43807878094991, 362.4692635312307)
)

$r = 2
foreach ($row in $rows) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Match the date-formatted style (column A) used on the other sheets
$wsWeekly.Range("A2:A50").Copy()
$wsForecast.Range("A2:A58").PasteSpecial(-4122)

# Restore the originally-active sheet (the new sheet's addition otherwise
# leaves itself activated/selected)
$wsWeekly.Activate()

